$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 558.3333
$ws.Range("I6").Value = 364.55554
$ws.Range("J6").Value = 1139.6666
$ws.Range("K6").Value = 1093.66662
$ws.Range("L6").Value = 3418.9998
$ws.Range("M6").Value = -981.66662
$ws.Range("N6").Value = -3642.9998
$ws.Range("H28").Value = 1748.4375
$ws.Range("I28").Value = 1269.5238
$ws.Range("J28").Value = 2662.7273
$ws.Range("K28").Value = 1269.5238
$ws.Range("L28").Value = 2662.7273
$ws.Range("M28").Value = -784.5237999999999
$ws.Range("N28").Value = -3632.7273
$ws.Range("H70").Value = 113876.22
$ws.Range("J70").Value = 203480.2
$ws.Range("L70").Value = 610440.6000000001
$ws.Range("N70").Value = -610980.6000000001
$ws.Range("H73").Value = 113876.22
$ws.Range("J73").Value = 203480.2
$ws.Range("L73").Value = 610440.6000000001
$ws.Range("N73").Value = -612312.6000000001
$ws.Range("H86").Value = 2644.818
$ws.Range("J86").Value = 2697.6
$ws.Range("L86").Value = 2697.6
$ws.Range("N86").Value = -4943.6
$ws.Range("H89").Value = 2644.818
$ws.Range("J89").Value = 2697.6
$ws.Range("L89").Value = 13488
$ws.Range("N89").Value = -24720
$ws.Range("H98").Value = 4133.35
$ws.Range("I98").Value = 1339.3334
$ws.Range("K98").Value = 1339.3334
$ws.Range("M98").Value = 158.6666
$ws.Range("H107").Value = 127.61539
$ws.Range("I107").Value = 125
$ws.Range("J107").Value = 136.33333
$ws.Range("K107").Value = 125
$ws.Range("L107").Value = 136.33333
$ws.Range("M107").Value = 1795
$ws.Range("N107").Value = -3976.33333
$ws.Range("H122").Value = 4133.35
$ws.Range("I122").Value = 1339.3334
$ws.Range("K122").Value = 4018.0002
$ws.Range("M122").Value = -1568.0002
$ws.Range("H135").Value = 1449.6
$ws.Range("I135").Value = 1708
$ws.Range("J135").Value = 846.6667
$ws.Range("K135").Value = 15372
$ws.Range("L135").Value = 7620.0003
$ws.Range("M135").Value = -12837
$ws.Range("N135").Value = -12690.0003
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7592.8237
$ws.Range("I32").Value = 886
$ws.Range("J32").Value = 29390
$ws.Range("K32").Value = 886
$ws.Range("L32").Value = 29390
$ws.Range("M32").Value = -599
$ws.Range("N32").Value = -29964
$ws.Range("H61").Value = 4647.6875
$ws.Range("I61").Value = 2985.2307
$ws.Range("K61").Value = 2985.2307
$ws.Range("M61").Value = -2773.2307
$ws.Range("H74").Value = 1739.5834
$ws.Range("I74").Value = 1471.9678
$ws.Range("J74").Value = 3398.8
$ws.Range("K74").Value = 1471.9678
$ws.Range("L74").Value = 3398.8
$ws.Range("M74").Value = -597.9677999999999
$ws.Range("N74").Value = -5146.8
$ws.Range("H77").Value = 1739.5834
$ws.Range("I77").Value = 1471.9678
$ws.Range("J77").Value = 3398.8
$ws.Range("K77").Value = 7359.839
$ws.Range("L77").Value = 16994
$ws.Range("M77").Value = -2991.839
$ws.Range("N77").Value = -25730
$ws.Range("H136").Value = 4647.6875
$ws.Range("I136").Value = 2985.2307
$ws.Range("K136").Value = 8955.6921
$ws.Range("M136").Value = -6405.6921
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6778.3228
$ws.Range("I31").Value = 6041.737
$ws.Range("J31").Value = 7944.5835
$ws.Range("K31").Value = 6041.737
$ws.Range("L31").Value = 7944.5835
$ws.Range("M31").Value = -5746.737
$ws.Range("N31").Value = -8534.583500000001
$ws.Range("H34").Value = 6778.3228
$ws.Range("I34").Value = 6041.737
$ws.Range("J34").Value = 7944.5835
$ws.Range("K34").Value = 6041.737
$ws.Range("L34").Value = 7944.5835
$ws.Range("M34").Value = -5839.737
$ws.Range("N34").Value = -8348.583500000001
$ws.Range("H134").Value = 1386.7084
$ws.Range("I134").Value = 1339.174
$ws.Range("J134").Value = 2480
$ws.Range("K134").Value = 4017.522
$ws.Range("L134").Value = 7440
$ws.Range("M134").Value = -1482.522
$ws.Range("N134").Value = -12510
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 45000
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 45000
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 45000
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -45302
$ws.Range("H46").Value = 21955.688
$ws.Range("I46").Value = 3715.1667
$ws.Range("J46").Value = 32900
$ws.Range("K46").Value = 3715.1667
$ws.Range("L46").Value = 32900
$ws.Range("M46").Value = -3559.1667
$ws.Range("N46").Value = -33212
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4774.853
$ws.Range("I46").Value = 5144
$ws.Range("J46").Value = 4641.96
$ws.Range("K46").Value = 5144
$ws.Range("L46").Value = 4641.96
$ws.Range("M46").Value = -4956
$ws.Range("N46").Value = -5017.96
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3399.6667
$ws.Range("J62").Value = 2899
$ws.Range("L62").Value = 2899
$ws.Range("N62").Value = -4147
$ws.Range("H65").Value = 3399.6667
$ws.Range("J65").Value = 2899
$ws.Range("L65").Value = 14495
$ws.Range("N65").Value = -20735
$ws.Range("H81").Value = 2350.2942
$ws.Range("I81").Value = 2064.0667
$ws.Range("J81").Value = 4497
$ws.Range("K81").Value = 4128.1334
$ws.Range("L81").Value = 8994
$ws.Range("M81").Value = -3067.1334
$ws.Range("N81").Value = -11116
$ws.Range("H84").Value = 2350.2942
$ws.Range("I84").Value = 2064.0667
$ws.Range("J84").Value = 4497
$ws.Range("K84").Value = 20640.667
$ws.Range("L84").Value = 44970
$ws.Range("M84").Value = -15336.667
$ws.Range("N84").Value = -55578
$ws.Range("H132").Value = 2788.7273
$ws.Range("I132").Value = 2570.1
$ws.Range("J132").Value = 4975
$ws.Range("K132").Value = 7710.299999999999
$ws.Range("L132").Value = 14925
$ws.Range("M132").Value = -5180.299999999999
$ws.Range("N132").Value = -19985
$ws.Range("H136").Value = 4990.2383
$ws.Range("I136").Value = 4877.5
$ws.Range("J136").Value = 5666.6665
$ws.Range("K136").Value = 14632.5
$ws.Range("L136").Value = 16999.9995
$ws.Range("M136").Value = -12082.5
$ws.Range("N136").Value = -22099.9995
